$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = 9767.7800000000007
$ws.Range("B9").Value = 9968.14
$ws.Range("C9").Value = 79.650000000000006
$ws.Range("D9").Value = 78.05
$ws.Range("E9").Value = $false
$ws.Range("F9").Value = -2.0099999999999998
$ws.Range("G9").Value = 42612.67291666667
$ws.Range("G9").NumberFormat = "m/d/yy h:mm"
$ws.Range("H9").Value = $false
